$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01114771857106678
$ws.Range("D2").Value = 0.111850647797489
$ws.Range("E2").Value = 0.09511822491632671
$ws.Range("F2").Value = 2.096748210335988
$ws.Range("G2").Value = 1.657605119417553
$ws.Range("H2").Value = 1.332276378139909
$ws.Range("I2").Value = 1.299107241881885
$ws.Range("J2").Value = 0.1070949409277731
$ws.Range("M2").Value = 1.591777731784489
$ws.Range("N2").Value = 1.662922876798859
$ws.Range("C3").Value = 0.01132455048978476
$ws.Range("D3").Value = 0.1098176463855651
$ws.Range("E3").Value = 0.09552677440884594
$ws.Range("F3").Value = 2.046759814905371
$ws.Range("G3").Value = 1.586021052199413
$ws.Range("H3").Value = 1.306690657001013
$ws.Range("I3").Value = 1.252525169104359
$ws.Range("J3").Value = 0.1090268536863483
$ws.Range("M3").Value = 1.446390440664544
$ws.Range("N3").Value = 1.541200146224611
$ws.Range("C4").Value = 0.01143788550734015
$ws.Range("D4").Value = 0.1085573265249309
$ws.Range("E4").Value = 0.09580432564145092
$ws.Range("F4").Value = 2.017739411520111
$ws.Range("G4").Value = 1.54346023831701
$ws.Range("H4").Value = 1.29195954484041
$ws.Range("I4").Value = 1.224970147524274
$ws.Range("J4").Value = 0.1102802625411172
$ws.Range("M4").Value = 1.357178730461627
$ws.Range("N4").Value = 1.466736462563745
$ws.Range("C5").Value = 0.01148526417181095
$ws.Range("D5").Value = 0.1080406870349222
$ws.Range("E5").Value = 0.09592417369745121
$ws.Range("F5").Value = 2.006330286910213
$ws.Range("G5").Value = 1.526462172412948
$ws.Range("H5").Value = 1.28620077282423
$ws.Range("I5").Value = 1.21400167476952
$ws.Range("J5").Value = 0.1108079189687441
$ws.Range("M5").Value = 1.320839995588926
$ws.Range("N5").Value = 1.436464336298684
$ws.Range("C6").Value = 0.01149320337312265
$ws.Range("D6").Value = 0.1079547147461639
$ws.Range("E6").Value = 0.09594448254688004
$ws.Range("F6").Value = 2.004460884946965
$ws.Range("G6").Value = 1.523660419617329
$ws.Range("H6").Value = 1.285259234413616
$ws.Range("I6").Value = 1.212196017529152
$ws.Range("J6").Value = 0.1108965550499104
$ws.Range("M6").Value = 1.314806973632358
$ws.Range("N6").Value = 1.431442148087768
$ws.Range("C7").Value = 0.01143851964269516
$ws.Range("D7").Value = 0.108550371314351
$ws.Range("E7").Value = 0.09580591460986554
$ws.Range("F7").Value = 2.017583860783034
$ws.Range("G7").Value = 1.543229601842143
$ws.Range("H7").Value = 1.291880893327061
$ws.Range("I7").Value = 1.224821172020555
$ws.Range("J7").Value = 0.1102873103674664
$ws.Range("M7").Value = 1.356688587725927
$ws.Range("N7").Value = 1.46632790428248
$ws.Range("C8").Value = 0.01120770051610143
$ws.Range("D8").Value = 0.1111521574500216
$ws.Range("E8").Value = 0.0952535689631322
$ws.Range("F8").Value = 2.079162798403829
$ws.Range("G8").Value = 1.632631411962052
$ws.Range("H8").Value = 1.323250201767564
$ws.Range("I8").Value = 1.282826851649176
$ws.Range("J8").Value = 0.1077471052990866
$ws.Range("M8").Value = 1.541637279139593
$ws.Range("N8").Value = 1.62089786779444
$ws.Range("C9").Value = 0.01079296309869071
$ws.Range("D9").Value = 0.116159749568709
$ws.Range("E9").Value = 0.094380972967814
$ws.Range("F9").Value = 2.2133673725412
$ws.Range("G9").Value = 1.819197354715897
$ws.Range("H9").Value = 1.392617601391748
$ws.Range("I9").Value = 1.405012008416207
$ws.Range("J9").Value = 0.1032997510086595
$ws.Range("M9").Value = 1.904724642456898
$ws.Range("N9").Value = 1.926056914493273
$ws.Range("C10").Value = 0.010511536349328
$ws.Range("D10").Value = 0.1197830795509134
$ws.Range("E10").Value = 0.09386648782318296
$ws.Range("F10").Value = 2.320425677567272
$ws.Range("G10").Value = 1.963425009782014
$ws.Range("H10").Value = 1.44849818184332
$ws.Range("I10").Value = 1.500119721073432
$ws.Range("J10").Value = 0.1003587224672131
$ws.Range("M10").Value = 2.171698886115252
$ws.Range("N10").Value = 2.151356550238233
$ws.Range("C11").Value = 0.0103886035004972
$ws.Range("D11").Value = 0.1214197510868189
$ws.Range("E11").Value = 0.09365955627412781
$ws.Range("F11").Value = 2.371023721356607
$ws.Range("G11").Value = 2.030658821529869
$ws.Range("H11").Value = 1.475016169076014
$ws.Range("I11").Value = 1.544589374823573
$ws.Range("J11").Value = 0.09909185871675064
$ws.Range("M11").Value = 2.293194683101632
$ws.Range("N11").Value = 2.254060850396343
$ws.Range("C12").Value = 0.01034278788222665
$ws.Range("D12").Value = 0.1220378764832333
$ws.Range("E12").Value = 0.0935850644097993
$ws.Range("F12").Value = 2.390461163480211
$ws.Range("G12").Value = 2.056357240924115
$ws.Range("H12").Value = 1.485217941091321
$ws.Range("I12").Value = 1.561605498555409
$ws.Range("J12").Value = 0.09862236489226284
$ws.Range("M12").Value = 2.339207965529511
$ws.Range("N12").Value = 2.292980275854802
$ws.Range("C13").Value = 0.01035262227048017
$ws.Range("D13").Value = 0.1219048250049894
$ws.Range("E13").Value = 0.09360093598262154
$ws.Range("F13").Value = 2.386262570573194
$ws.Range("G13").Value = 2.050811942941891
$ws.Range("H13").Value = 1.483013656901619
$ws.Range("I13").Value = 1.557932873489477
$ws.Range("J13").Value = 0.0987230229304572
$ws.Range("M13").Value = 2.329297961092607
$ws.Range("N13").Value = 2.284597108258311
$ws.Range("C14").Value = 0.01038481943843861
$ws.Range("D14").Value = 0.1214706375102992
$ws.Range("E14").Value = 0.0936533504215884
$ws.Range("F14").Value = 2.372617273691816
$ws.Range("G14").Value = 2.032768238940207
$ws.Range("H14").Value = 1.475852256202245
$ws.Range("I14").Value = 1.545985747099849
$ws.Range("J14").Value = 0.09905302779794312
$ws.Range("M14").Value = 2.296980122058102
$ws.Range("N14").Value = 2.257262241708361
$ws.Range("C15").Value = 0.01040463716788453
$ws.Range("D15").Value = 0.1212044713456422
$ws.Range("E15").Value = 0.09368595870229512
$ws.Range("F15").Value = 2.364295346369602
$ws.Range("G15").Value = 2.021747147664996
$ws.Range("H15").Value = 1.471486590789993
$ws.Range("I15").Value = 1.538690859649009
$ws.Range("J15").Value = 0.09925649965805228
$ws.Range("M15").Value = 2.277185182329532
$ws.Range("N15").Value = 2.240522336113429
$ws.Range("C16").Value = 0.01051967309759028
$ws.Range("D16").Value = 0.1196758876187474
$ws.Range("E16").Value = 0.09388055426827258
$ws.Range("F16").Value = 2.317157510151475
$ws.Range("G16").Value = 1.959064217956438
$ws.Range("H16").Value = 1.446787438832388
$ws.Range("I16").Value = 1.497238038087502
$ws.Range("J16").Value = 0.1004429464636125
$ws.Range("M16").Value = 2.163759690519441
$ws.Range("N16").Value = 2.144648613620348
$ws.Range("C17").Value = 0.01059155063594397
$ws.Range("D17").Value = 0.1187351951156259
$ws.Range("E17").Value = 0.09400685475148229
$ws.Range("F17").Value = 2.288728811185678
$ws.Range("G17").Value = 1.921029669399587
$ws.Range("H17").Value = 1.431917954203215
$ws.Range("I17").Value = 1.472118880993008
$ws.Range("J17").Value = 0.1011890031378542
$ws.Range("M17").Value = 2.094188057401453
$ws.Range("N17").Value = 2.085885733693772
$ws.Range("C18").Value = 0.01063337147065857
$ws.Range("D18").Value = 0.1181930405501106
$ws.Range("E18").Value = 0.09408205298637817
$ws.Range("F18").Value = 2.272555637536982
$ws.Range("G18").Value = 1.899305756041286
$ws.Range("H18").Value = 1.423468603308152
$ws.Range("I18").Value = 1.457784201588339
$ws.Range("J18").Value = 0.101624798774246
$ws.Range("M18").Value = 2.054176882074927
$ws.Range("N18").Value = 2.052107249104324
$ws.Range("C19").Value = 0.01064761338332021
$ws.Range("D19").Value = 0.1180092879466201
$ws.Range("E19").Value = 0.09410795329627319
$ws.Range("F19").Value = 2.267110165064366
$ws.Range("G19").Value = 1.891976464921328
$ws.Range("H19").Value = 1.420625455311722
$ws.Range("I19").Value = 1.452950070250708
$ws.Range("J19").Value = 0.101773498792376
$ws.Range("M19").Value = 2.04063061782324
$ws.Range("N19").Value = 2.040674034410131
$ws.Range("C20").Value = 0.01058384956482628
$ws.Range("D20").Value = 0.1188354464197374
$ws.Range("E20").Value = 0.09399314578390694
$ws.Range("F20").Value = 2.291736614219815
$ws.Range("G20").Value = 1.925062688025662
$ws.Range("H20").Value = 1.43349014171153
$ws.Range("I20").Value = 1.474781119855919
$ws.Range("J20").Value = 0.1011088922155139
$ws.Range("M20").Value = 2.1015936121577
$ws.Range("N20").Value = 2.092139058418127
$ws.Range("C21").Value = 0.01037534232742776
$ws.Range("D21").Value = 0.1215982134170233
$ws.Range("E21").Value = 0.09363785027479565
$ws.Range("F21").Value = 2.376617673991973
$ws.Range("G21").Value = 2.038061601580466
$ws.Range("H21").Value = 1.477951373894939
$ws.Range("I21").Value = 1.549490093660381
$ws.Range("J21").Value = 0.09895581924327423
$ws.Range("M21").Value = 2.306472517649496
$ws.Range("N21").Value = 2.265290433204825
$ws.Range("C22").Value = 0.01024336520469493
$ws.Range("D22").Value = 0.1233942542034185
$ws.Range("E22").Value = 0.09342817954737725
$ws.Range("F22").Value = 2.433709335764888
$ws.Range("G22").Value = 2.113305039855106
$ws.Range("H22").Value = 1.507942904560821
$ws.Range("I22").Value = 1.599346626302363
$ws.Range("J22").Value = 0.09760836083523827
$ws.Range("M22").Value = 2.44040393404876
$ws.Range("N22").Value = 2.378613884945992
$ws.Range("C23").Value = 0.01031340940808789
$ws.Range("D23").Value = 0.1224365434024932
$ws.Range("E23").Value = 0.09353803262724369
$ws.Range("F23").Value = 2.403089040415381
$ws.Range("G23").Value = 2.073017176584983
$ws.Range("H23").Value = 1.491849720508071
$ws.Range("I23").Value = 1.572641948923163
$ws.Range("J23").Value = 0.0983220532944955
$ws.Range("M23").Value = 2.368919817386796
$ws.Range("N23").Value = 2.318117546315762
$ws.Range("C24").Value = 0.01058732966964371
$ws.Range("D24").Value = 0.1187901269741261
$ws.Range("E24").Value = 0.09399933554988316
$ws.Range("F24").Value = 2.29037625459037
$ws.Range("G24").Value = 1.923238916335976
$ws.Range("H24").Value = 1.432779046560427
$ws.Range("I24").Value = 1.473577189689465
$ws.Range("J24").Value = 0.1011450889429248
$ws.Range("M24").Value = 2.098245602844656
$ws.Range("N24").Value = 2.08931191424324
$ws.Range("C25").Value = 0.01090108555619695
$ws.Range("D25").Value = 0.1148149899444206
$ws.Range("E25").Value = 0.09459467270800204
$ws.Range("F25").Value = 2.175595135691054
$ws.Range("G25").Value = 1.767491293152204
$ws.Range("H25").Value = 1.372998000595288
$ws.Range("I25").Value = 1.371034524048767
$ws.Range("J25").Value = 0.1044456194858459
$ws.Range("M25").Value = 1.806460185957036
$ws.Range("N25").Value = 1.843299734266282
